# Apply the update: remove the "于都·希佳微夏日文化交流会" event row from the
# "展览" (Exhibitions) sheet, and remove the corresponding
# "南昌·【8月24日】滑稽互动狂欢大作战《欢乐小丑嘉年华》" event row from the
# "全部类型" (All types) sheet. Deleting the entire row shifts all the rows
# below it up by one, which matches the canonical OOXML diff (dimension
# A1:I18 -> A1:I17 on the exhibitions sheet, and A1:I20 -> A1:I19 on the
# all-types sheet).

$wb = $excel.ActiveWorkbook

# "展览" sheet: data row 2 is "于都·希佳微夏日文化交流会" -> delete it.
$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Rows.Item(2).Delete()

# "全部类型" sheet: data row 3 is "南昌·【8月24日】滑稽互动狂欢大作战《欢乐小丑嘉年华》" -> delete it.
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Rows.Item(3).Delete()
